$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=11; I='ba'; J='Appreciation'},
    @{Row=26; I='ba'; J='Appreciation'},
    @{Row=29; I='sv'; J='Statement-opinion'},
    @{Row=38; I='ba'; J='Appreciation'},
    @{Row=41; I='aa'; J='Agree/Accept'},
    @{Row=58; I='ba'; J='Appreciation'},
    @{Row=59; I='sd'; J='Statement-non-opinion'},
    @{Row=64; I='sv'; J='Statement-opinion'},
    @{Row=76; I='ba'; J='Appreciation'},
    @{Row=78; I='sd'; J='Statement-non-opinion'},
    @{Row=84; I='aa'; J='Agree/Accept'},
    @{Row=87; I='sd'; J='Statement-non-opinion'},
    @{Row=90; I='ba'; J='Appreciation'},
    @{Row=97; I='ba'; J='Appreciation'},
    @{Row=100; I='ba'; J='Appreciation'},
    @{Row=110; I='sd'; J='Statement-non-opinion'},
    @{Row=113; I='sd'; J='Statement-non-opinion'},
    @{Row=120; I='sd'; J='Statement-non-opinion'},
    @{Row=131; I='aa'; J='Agree/Accept'},
    @{Row=142; I='sv'; J='Statement-opinion'},
    @{Row=143; I='sv'; J='Statement-opinion'},
    @{Row=158; I='aa'; J='Agree/Accept'},
    @{Row=193; I='sv'; J='Statement-opinion'},
    @{Row=208; I='sd'; J='Statement-non-opinion'},
    @{Row=218; I='ba'; J='Appreciation'},
    @{Row=256; I='sd'; J='Statement-non-opinion'},
    @{Row=257; I='sd'; J='Statement-non-opinion'},
    @{Row=267; I='sd'; J='Statement-non-opinion'},
    @{Row=278; I='sd'; J='Statement-non-opinion'},
    @{Row=314; I='sd'; J='Statement-non-opinion'},
    @{Row=323; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=324; I='aa'; J='Agree/Accept'},
    @{Row=353; I='sd'; J='Statement-non-opinion'},
    @{Row=362; I='ba'; J='Appreciation'},
    @{Row=378; I='ba'; J='Appreciation'},
    @{Row=385; I='sv'; J='Statement-opinion'},
    @{Row=387; I='aa'; J='Agree/Accept'},
    @{Row=395; I='sd'; J='Statement-non-opinion'},
    @{Row=412; I='ba'; J='Appreciation'},
    @{Row=420; I='sv'; J='Statement-opinion'},
    @{Row=423; I='aa'; J='Agree/Accept'},
    @{Row=437; I='aa'; J='Agree/Accept'},
    @{Row=439; I='aa'; J='Agree/Accept'},
    @{Row=453; I='sd'; J='Statement-non-opinion'},
    @{Row=462; I='ba'; J='Appreciation'},
    @{Row=484; I='aa'; J='Agree/Accept'},
    @{Row=485; I='%'; J='Uninterpretable'},
    @{Row=507; I='aa'; J='Agree/Accept'},
    @{Row=508; I='sd'; J='Statement-non-opinion'},
    @{Row=511; I='sv'; J='Statement-opinion'},
    @{Row=513; I='sv'; J='Statement-opinion'},
    @{Row=519; I='sv'; J='Statement-opinion'},
    @{Row=544; I='sv'; J='Statement-opinion'},
    @{Row=548; I='sv'; J='Statement-opinion'},
    @{Row=551; I='sd'; J='Statement-non-opinion'},
    @{Row=566; I='sv'; J='Statement-opinion'},
    @{Row=570; I='sd'; J='Statement-non-opinion'},
    @{Row=583; I='sv'; J='Statement-opinion'},
    @{Row=597; I='sv'; J='Statement-opinion'},
    @{Row=604; I='sv'; J='Statement-opinion'},
    @{Row=616; I='sv'; J='Statement-opinion'},
    @{Row=620; I='sd'; J='Statement-non-opinion'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

$wb.Save()
